# Update "想去人数" (interested-count) values in column F for the
# "展览" and "全部类型" sheets, reflecting a refreshed data scrape.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> [old value, new value] (old kept only for sanity-checking)
$updates = @{
    2  = 715
    3  = 63
    6  = 25
    7  = 27
    11 = 4605
    12 = 4419
    14 = 19
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
